$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 148
$ws.Range("D21").Value = 136
$ws.Range("E21").Value = 12
$ws.Range("F21").Value = 38.96848137535817
